$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "totaalstand" data (rows 2..21, columns A..J) after the upload.
# Columns: Rang, Speler, Score, 180'ers, 100+ finishes, Totaal Score, Aantal Darts,
#          3-Darts Gemiddelde, Totaal, Winnaar
$data = @(
    @(1,"Nathan May",40,3,2,23145,1199,57.91,45,0),
    @(2,"Yannick den Daggelder",34,1,0,23399,1337,52.5,35,2),
    @(3,"Nick Fitzpatrick",24,1,1,11558,544,63.74,26,1),
    @(4,"Alessandro Delia",22,0,0,15518,916,50.82,22,0),
    @(4,"Remco Riem",21,1,0,9836,469,62.92,22,1),
    @(6,"Louis Tweddle",17,0,1,11864,581,61.26,18,0),
    @(7,"Gijs Tromp",14,1,0,11451,624,55.05,15,0),
    @(8,"Rocky Van Den Eeckhoudt",12,2,0,7628,396,57.79,14,0),
    @(9,"Milan Schoenmakers",11,1,0,6984,333,62.92,12,0),
    @(10,"Max Walter",11,0,0,6262,477,39.38,11,0),
    @(10,"Afendi Kelana",8,1,2,8344,432,57.94,11,0),
    @(12,"Diego Meerveld",9,0,0,7022,457,46.1,9,0),
    @(12,"Francesco Petruccelli",8,1,0,6404,349,55.05,9,0),
    @(12,"martin Berry",8,1,0,4295,285,45.21,9,0),
    @(15,"Dartin Dan",6,1,0,5323,298,53.59,7,0),
    @(16,"Robin Willis",6,0,0,5327,323,49.48,6,0),
    @(17,"Magnus Gladh",4,0,0,4578,263,52.22,4,0),
    @(18,"Oscar Ebbeling",3,0,0,3730,258,43.37,3,0),
    @(18,"Tristan Snoep",3,0,0,4499,321,42.05,3,0),
    @(20,"Lee Dolan",1,0,0,1770,131,40.53,1,0)
)

$rowCount = $data.Count
$colCount = 10
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$ws.Range("A$startRow`:J$endRow").Value = $arr
